$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 512 (existing rows 512:535 shift down to 513:536)
$ws.Rows.Item(512).Insert()

# Populate the newly inserted row 512 with the new weekly Brócoli price record
$ws.Cells.Item(512, 1).Value = 5
$ws.Cells.Item(512, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(512, 3).Value = "Maule"
$ws.Cells.Item(512, 4).Value = 45041
$ws.Cells.Item(512, 5).Value = 7
$ws.Cells.Item(512, 6).Value = 100112023
$ws.Cells.Item(512, 7).Value = "Brócoli"
$ws.Cells.Item(512, 8).Value = "Sin especificar"
$ws.Cells.Item(512, 9).Value = "Primera"
$ws.Cells.Item(512, 10).Value = 5000
$ws.Cells.Item(512, 11).Value = 700
$ws.Cells.Item(512, 12).Value = 700
$ws.Cells.Item(512, 13).Value = 700
$ws.Cells.Item(512, 14).Value = "$/unidad"
$ws.Cells.Item(512, 15).Value = "Región del Maule"
$ws.Cells.Item(512, 16).Value = 700
$ws.Cells.Item(512, 17).Value = 1
$ws.Cells.Item(512, 18).Value = "Hortaliza"
